# Insert a new data row at row 50 (pushing existing rows 50..165 down to 51..166)
# and populate it with the new weekly price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("50:50").Insert()

$ws.Range("A50").Value = 9
$ws.Range("B50").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C50").Value = "Metropolitana"
$ws.Range("D50").Value = 44519
$ws.Range("E50").Value = 13
$ws.Range("F50").Value = 100112030
$ws.Range("G50").Value = "Poroto granado"
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 43
$ws.Range("K50").Value = 34000
$ws.Range("L50").Value = 36000
$ws.Range("M50").Value = 35023
$ws.Range("N50").Value = "`$/malla 25 kilos"
$ws.Range("O50").Value = "Perú"
$ws.Range("P50").Value = 1401
$ws.Range("Q50").Value = 25
$ws.Range("R50").Value = "Hortaliza"
